$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after "line6" (current row 7) / before "extr1" (current row 8),
# for the two new contingency lines "line7" and "line8". This pushes the existing
# extr1..extr8 rows down by two (from rows 8-15 to rows 10-17).
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()

# Copy the header-column formatting (bold + border, same as the other rows in column A)
# onto the two freshly inserted A cells so they match the rest of the table.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber / update the shifted extr1..extr8 rows (now rows 10-17).
# extr1 -> row 10 (values unchanged, only the running index moves)
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# extr2 -> row 11 (in_service flips False -> True)
$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# extr3 -> row 12 (unchanged values)
$ws.Range("A12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# extr4 -> row 13 (in_service flips False -> True)
$ws.Range("A13").Value = 11
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# extr5 -> row 14 (unchanged values)
$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# extr6 -> row 15 (unchanged values)
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# extr7 -> row 16 (in_service flips False -> True)
$ws.Range("A16").Value = 14
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# extr8 -> row 17 (unchanged values)
$ws.Range("A17").Value = 15
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
